$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 137, pushing the existing rows 137:224
# down to 139:226 (this also grows the sheet dimension to A1:R226).
$ws.Rows.Item(137).Insert()
$ws.Rows.Item(137).Insert()

# Populate the two newly inserted rows with the new weekly records.
# Row 137 - "Primera" quality record dated 2023-11-03 (serial 45233)
$ws.Range("A137").Value = 7
$ws.Range("B137").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C137").Value = "Ñuble"
$ws.Range("D137").Value = 45233
$ws.Range("E137").Value = 16
$ws.Range("F137").Value = 100112037
$ws.Range("G137").Value = "Cebollín"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 500
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 7000
$ws.Range("M137").Value = 6600
$ws.Range("N137").Value = "$/paquete 36 unidades"
$ws.Range("O137").Value = "Provincia de Diguillín"
$ws.Range("P137").Value = 183
$ws.Range("Q137").Value = 36
$ws.Range("R137").Value = "Hortaliza"

# Row 138 - "Segunda" quality record dated 2023-11-03 (serial 45233)
$ws.Range("A138").Value = 7
$ws.Range("B138").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C138").Value = "Ñuble"
$ws.Range("D138").Value = 45233
$ws.Range("E138").Value = 16
$ws.Range("F138").Value = 100112037
$ws.Range("G138").Value = "Cebollín"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Segunda"
$ws.Range("J138").Value = 300
$ws.Range("K138").Value = 5000
$ws.Range("L138").Value = 5000
$ws.Range("M138").Value = 5000
$ws.Range("N138").Value = "$/paquete 36 unidades"
$ws.Range("O138").Value = "Provincia de Diguillín"
$ws.Range("P138").Value = 139
$ws.Range("Q138").Value = 36
$ws.Range("R138").Value = "Hortaliza"
